$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# Paragraph index -> target slide index in the deck (chapter divider slides)
$targets = @{
    1 = 3    # Data_Type       -> Chapter 1
    2 = 25   # Control_Flow    -> Chapter 2
    3 = 33   # Processes       -> Chapter 3
    4 = 36   # Communication   -> Chapter 4
    5 = 43   # Interface       -> Chapter 5
    6 = 47   # Constraints     -> Chapter 6
    7 = 64   # Functional Coverage -> Chapter 7
}

foreach ($idx in 1..7) {
    $para = $tr.Paragraphs($idx, 1)
    $targetSlide = $p.Slides.Item($targets[$idx])
    $title = $targetSlide.Shapes.Title.TextFrame.TextRange.Text
    $actionSetting = $para.ActionSettings.Item(1)
    $hyperlink = $actionSetting.Hyperlink
    $hyperlink.Address = ""
    $hyperlink.SubAddress = $targetSlide.SlideID.ToString() + "," + $targetSlide.SlideIndex.ToString() + "," + $title
}
